$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column header "REX_DEF" in F1, matching the header style used by B1:E1
$ws.Range("F1").Value = "REX_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill F2:F19 with the default/empty-list marker "[]" for every data row
$ws.Range("F2:F19").Value = "[]"
